$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stat columns (B:H) for rows 2-47, leaving the Name column (A) and
# cell styles untouched.
$ws.Range("B2:H47").ClearContents()

# Match the author's final selection state.
$ws.Range("B2:H47").Select()
